# "updated G4 of Mavs Timberwolves game"
# Adds the Game-4 (Series 3) stat line for each player still missing it, and
# fixes Mike Conley's RegSeason J/K (row 269) which had been left at stale
# values before his own Game-4 row is appended.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# NOTE: this engine's PowerShell subset does not bind *named* parameters
# (e.g. "-r 85 -player 'X'") correctly - they come through empty - so this
# helper is always called positionally.
function Set-Row {
    param(
        $r,
        $player,
        $seriesNum,
        $game,
        $gameInSeries,
        $pts,
        $reb,
        $ast,
        $threePM,
        $ppg,
        $rpg,
        $apg
    )

    $ws.Cells.Item($r, 1).Value  = $player
    $ws.Cells.Item($r, 2).Value  = $seriesNum
    $ws.Cells.Item($r, 3).Value  = $game
    $ws.Cells.Item($r, 4).Value  = $gameInSeries
    $ws.Cells.Item($r, 5).Value  = $pts
    $ws.Cells.Item($r, 6).Value  = $reb
    $ws.Cells.Item($r, 7).Value  = $ast
    $ws.Cells.Item($r, 8).Value  = $threePM
    $ws.Cells.Item($r, 9).Value  = $ppg
    $ws.Cells.Item($r, 10).Value = $rpg
    $ws.Cells.Item($r, 11).Value = $apg
}

# New Game 4 (of Series 3) rows, one per player, inserted into the existing
# per-player gaps in the sheet.
Set-Row 85  "Luka Doncic"     16 3 4 28 15 10 4 33.9 9.2  9.8
Set-Row 104 "Kyrie Irving"    16 3 4 16 2  4  1 25.6 5    5.2
Set-Row 123 "PJ Washington"   16 3 4 10 5  2  2 12.9 5.6  1.9
Set-Row 142 "Dereck Lively"   16 3 4 "NA" "NA" "NA" "NA" 8.8 6.9 1.1
Set-Row 161 "Daniel Gafford"  16 3 4 12 8  1  0 11   7.6  1.6
Set-Row 180 "Josh Green"      16 3 4 5  1  1  1 8.2  3.2  2.3
Set-Row 198 "Anthony Edwards" 15 3 4 29 10 9  2 25.9 5.4  5.1
Set-Row 216 "KAT"             15 3 4 25 5  1  4 21.8 8.3  3
Set-Row 234 "Jaden McDaniels" 15 3 4 10 1  0  2 10.5 3.1  1.4
Set-Row 252 "Rudy Gobert"     15 3 4 13 10 1  0 14   12.9 1.3

# Mike Conley's RegSeason RPG/APG (cols J/K) on his existing last row were
# stale placeholders; correct them to match his RegSeasonPPG (col I) ...
$ws.Cells.Item(269, 10).Value = 11.4
$ws.Cells.Item(269, 11).Value = 11.4

# ... then append his own Game 4 (Series 3) row.
Set-Row 270 "Mike Conley" 15 3 4 14 3 7 1 11.4 11.4 11.4

# Match the reviewer's final on-screen selection/scroll position.
$ws.Range("I84:K85").Select()
